$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 29   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/24/2022  Through  10/30/2022"

# --- Crime data table updates (rows 16-27) ---
# Row 16
$ws.Range("I14").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 110
$ws.Range("J16").Value = 89
$ws.Range("K16").Value = 23.595505617977
$ws.Range("L16").Value = -17.910447761194
$ws.Range("M16").Value = -45.544554455445
$ws.Range("N16").Value = -81.132075471698

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = -17.647058823529
$ws.Range("I17").Value = 191
$ws.Range("J17").Value = 158
$ws.Range("K17").Value = 20.886075949367
$ws.Range("L17").Value = 34.507042253521
$ws.Range("M17").Value = 27.333333333333
$ws.Range("N17").Value = -67.125645438898

# Row 18
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -62.5
$ws.Range("I18").Value = 59
$ws.Range("J18").Value = 99
$ws.Range("K18").Value = -40.404040404040
$ws.Range("L18").Value = -25.316455696202
$ws.Range("M18").Value = -36.559139784946
$ws.Range("N18").Value = -90.937019969278

# Row 19
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 20
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = -28.571428571428
$ws.Range("I19").Value = 292
$ws.Range("J19").Value = 264
$ws.Range("K19").Value = 10.606060606060
$ws.Range("L19").Value = 32.126696832579
$ws.Range("M19").Value = 135.483870967742
$ws.Range("N19").Value = 3.180212014134

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 150
$ws.Range("I20").Value = 66
$ws.Range("J20").Value = 52
$ws.Range("K20").Value = 26.923076923076
$ws.Range("L20").Value = 43.478260869565
$ws.Range("M20").Value = 53.488372093023
$ws.Range("N20").Value = -75.373134328358

# Row 21
$ws.Range("C21").Value = 12
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 55
$ws.Range("G21").Value = 67
$ws.Range("H21").Value = -17.910447761194
$ws.Range("I21").Value = 738
$ws.Range("J21").Value = 669
$ws.Range("K21").Value = 10.313901345291
$ws.Range("L21").Value = 16.220472440944
$ws.Range("M21").Value = 15.132605304212
$ws.Range("N21").Value = -69.975589910496

# Row 22
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("L22").Value = 23.529411764705

# Row 23
$ws.Range("I14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = 1
$ws.Range("K14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100
$ws.Range("I14").Copy()
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("G23").Value = 1
$ws.Range("K14").Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("H23").Value = -100
$ws.Range("J23").Value = 9
$ws.Range("K23").Value = -33.333333333333
$ws.Range("L23").Value = -40

# Row 24
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -30.434782608695
$ws.Range("F24").Value = 83
$ws.Range("G24").Value = 90
$ws.Range("H24").Value = -7.777777777777
$ws.Range("I24").Value = 1001
$ws.Range("J24").Value = 725
$ws.Range("K24").Value = 38.068965517241
$ws.Range("L24").Value = 41.184767277856
$ws.Range("M24").Value = 202.416918429003

# Row 25
$ws.Range("C25").Value = 5
$ws.Range("E25").Value = -28.571428571428
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = -7.142857142857
$ws.Range("I25").Value = 298
$ws.Range("J25").Value = 270
$ws.Range("K25").Value = 10.370370370370
$ws.Range("L25").Value = 17.786561264822
$ws.Range("M25").Value = -30.046948356807

# Row 26
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 100
$ws.Range("L26").Value = 57.142857142857

# Row 27
$ws.Range("C27").Value = 4
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = 9
$ws.Range("H27").Value = 800
$ws.Range("I27").Value = 49
$ws.Range("K27").Value = 75
$ws.Range("L27").Value = 104.166666666667

$excel.CutCopyMode = 0
